$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.306.67"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.874.74"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7123"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08030"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3150"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08219"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "1.886.15"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.245"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.393"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008517"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.95%  "
$ws.Range("D18").Value = "29.316.15"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "2.132.83"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.764"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1558"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.037"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.398"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.301"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.175"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.939"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7633"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.176"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01875"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").Value = "1.265.23"
$ws.Range("E39").Value = "  +3.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.751"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.450"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9138"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "112.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("E45").Value = "  +9.81%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "2.031.08"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5231"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.798"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.472"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4351"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.00%  "
